$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 3.445154666666667
$ws.Cells.Item(2, 8).Value = 10.335464
$ws.Cells.Item(2, 9).Value = 0.01110365039942287
$ws.Cells.Item(2, 10).Value = 0.01110365039942286
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.724001666666667
$ws.Cells.Item(2, 14).Value = 8.172005
$ws.Cells.Item(2, 15).Value = 0.04635500474236593
$ws.Cells.Item(2, 16).Value = 0.04635500474236593
$ws.Cells.Item(2, 17).Value = 9.384607053924444
$ws.Cells.Item(2, 18).Value = 84.46146348532001
$ws.Cells.Item(2, 19).Value = 0.0005147097669228203
$ws.Cells.Item(2, 20).Value = 0.0005147097669228202

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 3.445154666666667
$ws.Cells.Item(3, 8).Value = 10.335464
$ws.Cells.Item(3, 9).Value = 0.01110365039942287
$ws.Cells.Item(3, 10).Value = 0.01110365039942286
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 40.62063066666667
$ws.Cells.Item(3, 14).Value = 121.861892
$ws.Cells.Item(3, 15).Value = 0.6912512390256352
$ws.Cells.Item(3, 16).Value = 0.6912512390256351
$ws.Cells.Item(3, 17).Value = 139.9443553042098
$ws.Cells.Item(3, 18).Value = 1259.499197737888
$ws.Cells.Item(3, 19).Value = 0.007675412096308546
$ws.Cells.Item(3, 20).Value = 0.007675412096308543

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 3.445154666666667
$ws.Cells.Item(4, 8).Value = 10.335464
$ws.Cells.Item(4, 9).Value = 0.01110365039942287
$ws.Cells.Item(4, 10).Value = 0.01110365039942286
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 15.419285
$ws.Cells.Item(4, 14).Value = 46.257855
$ws.Cells.Item(4, 15).Value = 0.2623937562319988
$ws.Cells.Item(4, 16).Value = 0.2623937562319988
$ws.Cells.Item(4, 17).Value = 53.12182167441334
$ws.Cells.Item(4, 18).Value = 478.09639506972
$ws.Cells.Item(4, 19).Value = 0.0029135285361915
$ws.Cells.Item(4, 20).Value = 0.002913528536191499

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 291.329961
$ws.Cells.Item(5, 8).Value = 873.989883
$ws.Cells.Item(5, 9).Value = 0.938949437922138
$ws.Cells.Item(5, 10).Value = 0.938949437922138
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.724001666666667
$ws.Cells.Item(5, 14).Value = 8.172005
$ws.Cells.Item(5, 15).Value = 0.04635500474236593
$ws.Cells.Item(5, 16).Value = 0.04635500474236593
$ws.Cells.Item(5, 17).Value = 793.583299313935
$ws.Cells.Item(5, 18).Value = 7142.249693825415
$ws.Cells.Item(5, 19).Value = 0.04352500564772253
$ws.Cells.Item(5, 20).Value = 0.04352500564772253

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 291.329961
$ws.Cells.Item(6, 8).Value = 873.989883
$ws.Cells.Item(6, 9).Value = 0.938949437922138
$ws.Cells.Item(6, 10).Value = 0.938949437922138
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 40.62063066666667
$ws.Cells.Item(6, 14).Value = 121.861892
$ws.Cells.Item(6, 15).Value = 0.6912512390256352
$ws.Cells.Item(6, 16).Value = 0.6912512390256351
$ws.Cells.Item(6, 17).Value = 11834.0067479154
$ws.Cells.Item(6, 18).Value = 106506.0607312386
$ws.Cells.Item(6, 19).Value = 0.6490499623461017
$ws.Cells.Item(6, 20).Value = 0.6490499623461016

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 291.329961
$ws.Cells.Item(7, 8).Value = 873.989883
$ws.Cells.Item(7, 9).Value = 0.938949437922138
$ws.Cells.Item(7, 10).Value = 0.938949437922138
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 15.419285
$ws.Cells.Item(7, 14).Value = 46.257855
$ws.Cells.Item(7, 15).Value = 0.2623937562319988
$ws.Cells.Item(7, 16).Value = 0.2623937562319988
$ws.Cells.Item(7, 17).Value = 4492.099697697885
$ws.Cells.Item(7, 18).Value = 40428.89727928097
$ws.Cells.Item(7, 19).Value = 0.2463744699283138
$ws.Cells.Item(7, 20).Value = 0.2463744699283138

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 15.497141
$ws.Cells.Item(8, 8).Value = 46.491423
$ws.Cells.Item(8, 9).Value = 0.04994691167843914
$ws.Cells.Item(8, 10).Value = 0.04994691167843914
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 2.724001666666667
$ws.Cells.Item(8, 14).Value = 8.172005
$ws.Cells.Item(8, 15).Value = 0.04635500474236593
$ws.Cells.Item(8, 16).Value = 0.04635500474236593
$ws.Cells.Item(8, 17).Value = 42.21423791256833
$ws.Cells.Item(8, 18).Value = 379.928141213115
$ws.Cells.Item(8, 19).Value = 0.002315289327720579
$ws.Cells.Item(8, 20).Value = 0.002315289327720579

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 15.497141
$ws.Cells.Item(9, 8).Value = 46.491423
$ws.Cells.Item(9, 9).Value = 0.04994691167843914
$ws.Cells.Item(9, 10).Value = 0.04994691167843914
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 40.62063066666667
$ws.Cells.Item(9, 14).Value = 121.861892
$ws.Cells.Item(9, 15).Value = 0.6912512390256352
$ws.Cells.Item(9, 16).Value = 0.6912512390256351
$ws.Cells.Item(9, 17).Value = 629.5036409502574
$ws.Cells.Item(9, 18).Value = 5665.532768552316
$ws.Cells.Item(9, 19).Value = 0.03452586458322503
$ws.Cells.Item(9, 20).Value = 0.03452586458322502

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 15.497141
$ws.Cells.Item(10, 8).Value = 46.491423
$ws.Cells.Item(10, 9).Value = 0.04994691167843914
$ws.Cells.Item(10, 10).Value = 0.04994691167843914
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 15.419285
$ws.Cells.Item(10, 14).Value = 46.257855
$ws.Cells.Item(10, 15).Value = 0.2623937562319988
$ws.Cells.Item(10, 16).Value = 0.2623937562319988
$ws.Cells.Item(10, 17).Value = 238.954833764185
$ws.Cells.Item(10, 18).Value = 2150.593503877665
$ws.Cells.Item(10, 19).Value = 0.01310575776749354
$ws.Cells.Item(10, 20).Value = 0.01310575776749354
